# Refresh the cryptos list: updated prices / 1h volume % for most rows,
# plus a few coins that moved around in the ranking (rows 30, 31, 51).
#
# Columns D (Price) and E (Volume(1h)) are stored as plain text in this
# sheet (e.g. "71.178.96", "  +1.09%  "), not numbers - some of the new
# Price values look like ordinary decimals (e.g. "694.10"), and Excel's
# COM layer auto-converts a bare numeric-looking string typed into
# .Value into a real number. To keep those particular cells as text
# (matching the sheet's existing convention) we momentarily force the
# cell to Text format before assigning, then restore the cell's style
# from its neighboring (always-text) B-column cell so no stray number
# format sticks around on the cell itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.178.96"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "3.860.91"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "694.10"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.41"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("D7").Value = "3.858.03"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.37"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +6.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.66"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("D15").Value = "4.515.28"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "3.862.73"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "71.249.15"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.80"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  -3.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "495.00"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  +3.65%  "
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.84"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.36"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.63"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  +3.16%  "
$ws.Range("E28").Value = "  +1.91%  "
$ws.Range("D29").Value = "4.019.19"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.15"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +10.70%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.65"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  +3.46%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.76"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.183"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.32"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("D37").Value = "3.815.58"
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("E40").Value = "  +13.32%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.07"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  +1.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.02"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  +5.82%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "164.76"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000308"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  +4.93%  "
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.303"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "418.71"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  +6.71%  "
